$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.615.85"
$ws.Range("E2").Value = "  +1.83%  "
$ws.Range("D3").Value = "3.265.08"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'583.40"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").Value = "'181.40"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").Value = "'0.133"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").Value = "'6.66"
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "'0.424"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").Value = "3.834.74"
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("E13").Value = "  +0.26%  "
$ws.Range("D14").Value = "'28.47"
$ws.Range("E14").Value = "  +0.11%  "
$ws.Range("D15").Value = "68.605.29"
$ws.Range("E15").Value = "  +1.85%  "
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").Value = "3.266.33"
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("D19").Value = "'13.56"
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").Value = "'393.99"
$ws.Range("E20").Value = "  +4.87%  "
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("D22").Value = "'72.05"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("D24").Value = "'0.516"
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  +4.48%  "
$ws.Range("D27").Value = "'9.60"
$ws.Range("E27").Value = "  +0.38%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").Value = "'1.97"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").Value = "'5.69"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").Value = "'22.94"
$ws.Range("E31").Value = "  +1.24%  "
$ws.Range("D32").Value = "'7.13"
$ws.Range("E32").Value = "  +3.13%  "
$ws.Range("E33").Value = "  +0.72%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'164.04"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  +3.64%  "
$ws.Range("D38").Value = "'0.829"
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").Value = "'26.42"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'4.60"
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("D41").Value = "'6.56"
$ws.Range("E41").Value = "  -3.80%  "
$ws.Range("E42").Value = "  -4.23%  "
$ws.Range("D43").Value = "'348.45"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").Value = "'41.33"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").Value = "2.610.02"
$ws.Range("E46").Value = "  -3.91%  "
$ws.Range("D47").Value = "'24.61"
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("D48").Value = "'0.0282"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "'6.32"
$ws.Range("E49").Value = "  +2.78%  "
$ws.Range("D50").Value = "'31.60"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("E51").Value = "  -0.17%  "
